# Einarbeitung der neuen Diagramme
# Adds row 9 with "reference value" formulas derived from row 7 (the
# column means), rescaled to a common basis (divide by 1.45, multiply by 5),
# and moves the active selection to H14 as left by the author after
# inserting the new charts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9 holds its own (non-shared) formula...
$ws.Range("B9").Formula = "=B`$7 / 1.45 * 5"

# ...while C9:I9 share one formula definition (C9 is the master cell).
$ws.Range("C9:I9").Formula = "=C`$7 / 1.45 * 5"

# Restore the cursor position recorded in the saved workbook.
[void]$ws.Range("H14").Select()
